$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 409, shifting existing rows 409:473 down to 410:474
$ws.Rows.Item(409).Insert()

# Populate the new row 409 - it shares most categorical attributes with the
# row right below it (old row 409, now shifted to row 410), but carries a
# newer date/price observation.
$ws.Cells.Item(409, 1).Value = 5                                         # A: Mercado ID
$ws.Cells.Item(409, 2).Value = "Macroferia Regional de Talca"            # B: Mercado
$ws.Cells.Item(409, 3).Value = "Maule"                                   # C: Region
$ws.Cells.Item(409, 4).Value = 45218                                     # D: Fecha
$ws.Cells.Item(409, 4).NumberFormat = $ws.Cells.Item(410, 4).NumberFormat
$ws.Cells.Item(409, 5).Value = 7                                         # E: Codreg
$ws.Cells.Item(409, 6).Value = "Fruta"                                   # F: Tipo
$ws.Cells.Item(409, 7).Value = 100108                                    # G: Producto ID
$ws.Cells.Item(409, 8).Value = "Tropicales y subtropicales"              # H: Producto
$ws.Cells.Item(409, 9).Value = 100108005                                 # I: Categoria ID
$ws.Cells.Item(409, 10).Value = "Piña"                                   # J: Categoria
$ws.Cells.Item(409, 11).Value = "Caramelo"                               # K: Variedad
$ws.Cells.Item(409, 12).Value = "Segunda"                                # L: Calidad
$ws.Cells.Item(409, 13).Value = 200                                      # M: Volumen
$ws.Cells.Item(409, 14).Value = 21000                                    # N: Precio minimo
$ws.Cells.Item(409, 15).Value = 21000                                    # O: Precio maximo
$ws.Cells.Item(409, 16).Value = 21000                                    # P: Precio promedio ponderado
$ws.Cells.Item(409, 17).Value = "$/caja 14 unidades"                     # Q: Unidad de comercializacion
$ws.Cells.Item(409, 18).Value = "Ecuador"                                # R: Origen
$ws.Cells.Item(409, 19).Value = 1500                                     # S: Precio $/Kg
$ws.Cells.Item(409, 20).Value = 14                                       # T: Kg / unidad
